# Applies the Spanish translations described by the commit:
#   Welcome              -> Bienvenida
#   What is SmartCash?   -> ¿Qué es SmartCash?
#   Brochure/Whitepaper  -> Folleto/Papel blanco
#   Hive Teams           -> Equipos de Hive
#   Roadmap              -> El mapa vial
#
# A plain Find/Replace (Range.Text = "...") works for most of these, but it
# re-serializes the owning run's text node from scratch, which has two side
# effects we must avoid:
#   1. If the run begins with a separate <w:tab/> element (as in the
#      "Brochure/Whitepaper" run), replacing just the trailing text merges
#      that tab into the new <w:t> as a literal "\t" character, losing the
#      dedicated <w:tab/> element.
#   2. xml:space="preserve" is only emitted by plain replacement when the
#      new string itself has leading/trailing whitespace, whereas the
#      target XML keeps/adds it on several runs regardless.
#
# To reproduce the exact target markup we instead rebuild the affected
# run(s) via Range.InsertXML, explicitly emitting the <w:tab/> element and
# the xml:space attribute exactly as required.

$d = $word.ActiveDocument

function New-RunPackageXml($innerRunBodyXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p><w:r>' + $innerRunBodyXml + '</w:r></w:p></w:body></w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'
}

function Replace-RunText($doc, $findText, $newText, $rPrXml, $preserve) {
    $found = $doc.Content
    $found.Find.Execute($findText)
    # Re-wrap the located bounds in a fresh Range; InsertXML behaves
    # unreliably when invoked directly on the Range returned by Find.
    $target = $doc.Range($found.Start, $found.End)

    if ($preserve) {
        $tXml = '<w:t xml:space="preserve">' + $newText + '</w:t>'
    } else {
        $tXml = '<w:t>' + $newText + '</w:t>'
    }

    $target.InsertXML((New-RunPackageXml ($rPrXml + $tXml)))
}

# Common run formatting shared by the "Welcome", "What is SmartCash?",
# "Hive Teams" and "Roadmap" runs.
$rPrArial = '<w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" w:eastAsia="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'

Replace-RunText $d "Welcome" "Bienvenida" $rPrArial $false
Replace-RunText $d "What is SmartCash?" "¿Qué es SmartCash?" $rPrArial $true
Replace-RunText $d "Hive Teams" "Equipos de Hive" $rPrArial $true
Replace-RunText $d "Roadmap" "El mapa vial" $rPrArial $true

# "Brochure/Whitepaper" lives in a run that also carries a leading
# <w:tab/> element. Locate the whole paragraph (tab + text, excluding the
# trailing paragraph mark) so the rebuilt run keeps the tab as a distinct
# element instead of collapsing it into the text.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Brochure/Whitepaper*") {
        $brochurePara = $p
        break
    }
}

$brochureRange = $d.Range($brochurePara.Range.Start, $brochurePara.Range.End - 1)
$rPrTimes = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:eastAsia="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'
$brochureRunXml = $rPrTimes + '<w:tab/><w:t xml:space="preserve">Folleto/Papel blanco</w:t>'
$brochureRange.InsertXML((New-RunPackageXml $brochureRunXml))
